# Generate Report for Archive
#
# The localization status text moves from "Ready for handoff" to
# "In Translation" everywhere it appears (Overview!E2:F4, zh-cn!C2:C4,
# de-de!C2:C4). Because the new text is shorter, the status/language
# columns that held it are re-narrowed to fit.
#
# Note: ColumnWidth is quantized by Excel to a whole-pixel character grid,
# so 12.5 is the closest settable value that reproduces the narrower
# target column width used by the generated report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn) and F (de-de), rows 2-4 ---
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2:F4").Value = "In Translation"
$ovw.Range("E1:F1").ColumnWidth = 12.5

# --- zh-cn sheet: column C (Status), rows 2-4 ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2:C4").Value = "In Translation"
$zh.Range("C1").ColumnWidth = 12.5

# --- de-de sheet: column C (Status), rows 2-4 ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2:C4").Value = "In Translation"
$de.Range("C1").ColumnWidth = 12.5
